$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Move the "_GoBack" bookmark from the end of the "Once the form is
#    submitted..." paragraph to the middle of the "I added some information
#    about myself..." paragraph, splitting the text right before "picture".
# ---------------------------------------------------------------------------
$oldBookmark = $d.Bookmarks.Item("_GoBack")
$oldBookmark.Delete()

$profilePara = $d.Paragraphs.Item(21)
$pictureRange = $profilePara.Range.Duplicate()
$pictureRange.Find.Execute("picture", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$bookmarkPoint = $d.Range($pictureRange.Start, $pictureRange.Start)
$d.Bookmarks.Add("_GoBack", $bookmarkPoint) | Out-Null

# ---------------------------------------------------------------------------
# 2) After the "All images were added into the images folder" paragraph,
#    append: a blank paragraph, a "Github Repo: <link>" paragraph, and
#    another blank paragraph.
# ---------------------------------------------------------------------------
$imagesPara = $d.Paragraphs.Item(22)

# Blank paragraph right after the images paragraph.
$imagesPara.Range.InsertParagraphAfter()
$blankPara1 = $d.Paragraphs.Item(23)
$blankPara1.Range.ListFormat.RemoveNumbers()
$blankPara1.Style = "Normal"

# Paragraph that will hold the "Github Repo: " text + hyperlink.
$blankPara1.Range.InsertParagraphAfter()
$repoPara = $d.Paragraphs.Item(24)
$repoPara.Range.ListFormat.RemoveNumbers()
$repoPara.Style = "Normal"

$repoUrl = "https://github.com/pulakazad/Blog_Project"

$repoInsert = $repoPara.Range.Duplicate()
$repoInsert.Collapse(1)
$repoInsert.InsertAfter("Github Repo: " + $repoUrl)

$repoLinkRange = $repoPara.Range.Duplicate()
$repoLinkRange.Find.Execute($repoUrl, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$d.Hyperlinks.Add($repoLinkRange, $repoUrl, "", "", $repoUrl) | Out-Null

# Trailing blank paragraph after the repo paragraph.
$repoPara.Range.InsertParagraphAfter()
$blankPara2 = $d.Paragraphs.Item(25)
$blankPara2.Range.ListFormat.RemoveNumbers()
$blankPara2.Style = "Normal"

# ---------------------------------------------------------------------------
# 3) Register the "Hyperlink" (and latent "Unresolved Mention") character
#    styles that Word normally pulls in automatically once a hyperlink is
#    present in the document.
# ---------------------------------------------------------------------------
$hyperlinkStyle = $d.Styles.Add("Hyperlink", 2)
$hyperlinkStyle.BaseStyle = "DefaultParagraphFont"
$hyperlinkStyle.Priority = 99
$hyperlinkStyle.UnhideWhenUsed = $true
$hyperlinkStyle.Font.Color = 0xC16305
$hyperlinkStyle.Font.Underline = 1

$mentionStyle = $d.Styles.Add("Unresolved Mention", 2)
$mentionStyle.BaseStyle = "DefaultParagraphFont"
$mentionStyle.Priority = 99
$mentionStyle.UnhideWhenUsed = $true
$mentionStyle.Hidden = $true
$mentionStyle.Font.Color = 0x605e5c

Write-Output "edit applied"
